$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.188575983047485
$ws.Range("B1").Value = 2.297712802886963
$ws.Range("C1").Value = 2.34470009803772
$ws.Range("D1").Value = 3.110618591308594
$ws.Range("E1").Value = 2.653542041778564
